# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2304"
#   "<header>_new" -> "<header>_FV2310"
# then turn the data range into a native Excel Table (with AutoFilter) and
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the header-row cells: swap the "_old"/"_new" suffixes for the
#    matching AHB format-version tags. "diff" (column K) is untouched.
# ---------------------------------------------------------------------
$lastCol = $ws.UsedRange.Columns.Count
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val -like "*_old") {
            $cell.Value = ($val -replace "_old$", "_FV2304")
        } elseif ($val -like "*_new") {
            $cell.Value = ($val -replace "_new$", "_FV2310")
        }
    }
}

# ---------------------------------------------------------------------
# 2. Freeze the header row (split below row 1).
# ---------------------------------------------------------------------
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 3. Turn the used range into an Excel Table with an AutoFilter, using
#    the (now renamed) header row as the column headers, and drop the
#    automatically-assigned table style so the data keeps its existing
#    cell formatting.
# ---------------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

[void]$ws.Range("A1").Select()
